$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11:21 down to 12:22.
$ws.Rows("11").Insert()

# Copy the style (number format) from the date cell that was just pushed down (now D12)
# onto the new D11 cell so the date formatting is preserved.
$ws.Range("D12").Copy()
$ws.Range("D11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 11 with the new weekly record.
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11, 3).Value = "Bíobío"
$ws.Cells.Item(11, 4).Value = 44488
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = "Jengibre"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 40
$ws.Cells.Item(11, 11).Value = 16000
$ws.Cells.Item(11, 12).Value = 17000
$ws.Cells.Item(11, 13).Value = 16500
$ws.Cells.Item(11, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 1269
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = "Hortaliza"
